{"js": "// Fix the project title typo: \"ArmagHEADon\" -> \"ArmagHEADdon\" (missing \"d\").\nconst results = context.document.body.search(\"ArmagHEADon\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Text \"ArmagHEADon\" not found in document body.');\n}\n\nconst target = results.items[0];\ntarget.insertText(\"ArmagHEADdon\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix the project title typo: \"ArmagHEADon\" -> \"ArmagHEADdon\" (missing \"d\").\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"ArmagHEADon\"\n$find.Replacement.Text = \"ArmagHEADdon\"\n$find.Forward = $true\n$find.Wrap = $wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n"}
